# Workbook / sheet references
$wb = $excel.ActiveWorkbook
$wsFin = $wb.Worksheets.Item(1)
$wsFin.Name = "TestFinancials"

# --- Update TestFinancials: add columns E (Updated) and F (2026E) ---
$wsFin.Range("E1").Value = "Updated"
$wsFin.Range("F1").Value = "2026E"

$wsFin.Range("E2").Value = 1600000
$wsFin.Range("F2").Value = 1800000

# --- Add BalanceSheet worksheet right after TestFinancials ---
$wsBS = $wb.Worksheets.Add($null, $wsFin)
$wsBS.Name = "BalanceSheet"

# Header years look numeric, force them to be stored as text like the source sheet
$wsBS.Range("B1:C1").NumberFormat = "@"

$wsBS.Range("A1").Value = "Item"
$wsBS.Range("B1").Value = "2023"
$wsBS.Range("C1").Value = "2024"

$wsBS.Range("A2").Value = "Cash"
$wsBS.Range("B2").Value = 50000
$wsBS.Range("C2").Value = 75000

$wsBS.Range("A3").Value = "Accounts Receivable"
$wsBS.Range("B3").Value = 100000
$wsBS.Range("C3").Value = 120000

$wsBS.Range("A4").Value = "Total Assets"
$wsBS.Range("B4").Value = 500000
$wsBS.Range("C4").Value = 600000

$wsBS.Range("A5").Value = "Debt"
$wsBS.Range("B5").Value = 200000
$wsBS.Range("C5").Value = 250000

$wsBS.Range("A6").Value = "Equity"
$wsBS.Range("B6").Value = 300000
$wsBS.Range("C6").Value = 350000

# --- Add CashFlow worksheet right after BalanceSheet ---
$wsCF = $wb.Worksheets.Add($null, $wsBS)
$wsCF.Name = "CashFlow"

$wsCF.Range("B1:C1").NumberFormat = "@"

$wsCF.Range("A1").Value = "Item"
$wsCF.Range("B1").Value = "2023"
$wsCF.Range("C1").Value = "2024"

$wsCF.Range("A2").Value = "Operating CF"
$wsCF.Range("B2").Value = 120000
$wsCF.Range("C2").Value = 150000

$wsCF.Range("A3").Value = "Investing CF"
$wsCF.Range("B3").Value = -80000
$wsCF.Range("C3").Value = -100000

$wsCF.Range("A4").Value = "Financing CF"
$wsCF.Range("B4").Value = -20000
$wsCF.Range("C4").Value = -30000

$wsCF.Range("A5").Value = "Net Change in Cash"
$wsCF.Range("B5").Value = 20000
$wsCF.Range("C5").Value = 20000
